$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 164.75
$ws.Range("I9").Value = 193
$ws.Range("K9").Value = 193
$ws.Range("M9").Value = -24
$ws.Range("H26").Value = 12500
$ws.Range("J26").Value = 12500
$ws.Range("L26").Value = 12500
$ws.Range("N26").Value = -13188
$ws.Range("H28").Value = 679.2143
$ws.Range("I28").Value = 332.5
$ws.Range("J28").Value = 1025.9286
$ws.Range("K28").Value = 332.5
$ws.Range("L28").Value = 1025.9286
$ws.Range("M28").Value = 152.5
$ws.Range("N28").Value = -1995.9286
$ws.Range("H40").Value = 45616.26
$ws.Range("I40").Value = 85316.586
$ws.Range("J40").Value = 2306.818
$ws.Range("K40").Value = 85316.586
$ws.Range("L40").Value = 2306.818
$ws.Range("M40").Value = -85141.586
$ws.Range("N40").Value = -2656.818
$ws.Range("H74").Value = 2959.8
$ws.Range("I74").Value = 3025
$ws.Range("K74").Value = 3025
$ws.Range("M74").Value = -2089
$ws.Range("H77").Value = 2959.8
$ws.Range("I77").Value = 3025
$ws.Range("K77").Value = 15125
$ws.Range("M77").Value = -10445
$ws.Range("H112").Value = 1035.1428
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 1035.1428
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 3105.4284
$ws.Range("M112").ClearContents()
$ws.Range("N112").Value = -5321.428400000001
$ws.Range("H131").Value = 4290.4043
$ws.Range("I131").Value = 1105
$ws.Range("K131").Value = 3315
$ws.Range("M131").Value = 1725

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 83.75
$ws.Range("I5").Value = 78
$ws.Range("J5").Value = 93.333336
$ws.Range("K5").Value = 78
$ws.Range("L5").Value = 93.333336
$ws.Range("M5").Value = 34
$ws.Range("N5").Value = -317.333336
$ws.Range("H22").Value = 425.14285
$ws.Range("I22").Value = 425.14285
$ws.Range("K22").Value = 425.14285
$ws.Range("M22").Value = -126.14285
$ws.Range("H25").Value = 10000
$ws.Range("I25").Value = 10000
$ws.Range("K25").Value = 10000
$ws.Range("M25").Value = -9598
$ws.Range("H33").Value = 22000
$ws.Range("I33").Value = 20000
$ws.Range("J33").Value = 24000
$ws.Range("K33").Value = 20000
$ws.Range("L33").Value = 24000
$ws.Range("M33").Value = -19671
$ws.Range("N33").Value = -24658
$ws.Range("H97").Value = 54214.58
$ws.Range("I97").Value = 63428.438
$ws.Range("J97").Value = 5074
$ws.Range("K97").Value = 63428.438
$ws.Range("L97").Value = 5074
$ws.Range("M97").Value = -62932.438
$ws.Range("N97").Value = -6066
$ws.Range("H132").Value = 1682.7778
$ws.Range("I132").Value = 1524.806
$ws.Range("J132").Value = 3799.6
$ws.Range("K132").Value = 4574.418
$ws.Range("L132").Value = 11398.8
$ws.Range("M132").Value = -2044.418
$ws.Range("N132").Value = -16458.8

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 83.75
$ws.Range("I4").Value = 78
$ws.Range("J4").Value = 93.333336
$ws.Range("K4").Value = 78
$ws.Range("L4").Value = 93.333336
$ws.Range("M4").Value = 37
$ws.Range("N4").Value = -323.333336
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()
$ws.Range("H24").Value = 190
$ws.Range("I24").Value = 190
$ws.Range("K24").Value = 190
$ws.Range("M24").Value = 45
$ws.Range("H29").Value = 177.33333
$ws.Range("I29").Value = 177.33333
$ws.Range("K29").Value = 177.33333
$ws.Range("M29").Value = 111.66667
$ws.Range("H33").Value = 1677.4286
$ws.Range("I33").Value = 873.6667
$ws.Range("K33").Value = 873.6667
$ws.Range("M33").Value = -537.6667

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()
$ws.Range("H31").Value = 27589
$ws.Range("I31").Value = 1068.762
$ws.Range("J31").Value = 41512.125
$ws.Range("K31").Value = 1068.762
$ws.Range("L31").Value = 41512.125
$ws.Range("M31").Value = -773.7619999999999
$ws.Range("N31").Value = -42102.125
$ws.Range("H34").Value = 27589
$ws.Range("I34").Value = 1068.762
$ws.Range("J34").Value = 41512.125
$ws.Range("K34").Value = 1068.762
$ws.Range("L34").Value = 41512.125
$ws.Range("M34").Value = -866.7619999999999
$ws.Range("N34").Value = -41916.125
$ws.Range("H58").Value = 3873.6792
$ws.Range("I58").Value = 961.44183
$ws.Range("J58").Value = 16396.3
$ws.Range("K58").Value = 961.44183
$ws.Range("L58").Value = 16396.3
$ws.Range("M58").Value = -758.44183
$ws.Range("N58").Value = -16802.3
$ws.Range("H122").Value = 782.25
$ws.Range("I122").Value = 759.3333
$ws.Range("J122").Value = 851
$ws.Range("K122").Value = 2277.9999
$ws.Range("L122").Value = 2553
$ws.Range("M122").Value = 172.0001000000002
$ws.Range("N122").Value = -7453
$ws.Range("H136").Value = 3873.6792
$ws.Range("I136").Value = 961.44183
$ws.Range("J136").Value = 16396.3
$ws.Range("K136").Value = 2884.32549
$ws.Range("L136").Value = 49188.89999999999
$ws.Range("M136").Value = -334.3254900000002
$ws.Range("N136").Value = -54288.89999999999

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 352.55554
$ws.Range("I18").Value = 319.53333
$ws.Range("J18").Value = 517.6667
$ws.Range("K18").Value = 958.5999899999999
$ws.Range("L18").Value = 1553.0001
$ws.Range("M18").Value = -789.5999899999999
$ws.Range("N18").Value = -1891.0001
$ws.Range("H69").Value = 3000
$ws.Range("J69").Value = 3000
$ws.Range("L69").Value = 9000
$ws.Range("N69").Value = -10622
$ws.Range("H72").Value = 3000
$ws.Range("J72").Value = 3000
$ws.Range("L72").Value = 27000
$ws.Range("N72").Value = -35112
$ws.Range("H111").Value = 2621.2856
$ws.Range("I111").Value = 1983
$ws.Range("J111").Value = 3100
$ws.Range("K111").Value = 5949
$ws.Range("L111").Value = 9300
$ws.Range("M111").Value = -2882
$ws.Range("N111").Value = -15434
$ws.Range("H112").Value = 144928.58
$ws.Range("I112").Value = 501000
$ws.Range("J112").Value = 2500
$ws.Range("K112").Value = 1503000
$ws.Range("L112").Value = 7500
$ws.Range("M112").Value = -1501892
$ws.Range("N112").Value = -9716
$ws.Range("H116").Value = 1484.875
$ws.Range("I116").Value = 813.1667
$ws.Range("J116").Value = 3500
$ws.Range("K116").Value = 2439.5001
$ws.Range("L116").Value = 10500
$ws.Range("M116").Value = 1002.4999
$ws.Range("N116").Value = -17384
$ws.Range("H122").Value = 11742
$ws.Range("I122").Value = 340
$ws.Range("K122").Value = 3060
$ws.Range("M122").Value = -610
$ws.Range("H131").Value = 872.01
$ws.Range("I131").Value = 730
$ws.Range("J131").Value = 873.44446
$ws.Range("K131").Value = 2190
$ws.Range("L131").Value = 2620.33338
$ws.Range("M131").Value = 2850
$ws.Range("N131").Value = -12700.33338

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H100").Value = 35886.668
$ws.Range("J100").Value = 35886.668
$ws.Range("L100").Value = 35886.668
$ws.Range("N100").Value = -38050.668
$ws.Range("H132").Value = 2116.1052
$ws.Range("I132").Value = 1290
$ws.Range("J132").Value = 3906
$ws.Range("K132").Value = 3870
$ws.Range("L132").Value = 11718
$ws.Range("M132").Value = -1340
$ws.Range("N132").Value = -16778

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5261.3
$ws.Range("I7").Value = 2000
$ws.Range("J7").Value = 6659
$ws.Range("K7").Value = 2000
$ws.Range("L7").Value = 6659
$ws.Range("M7").Value = -1888
$ws.Range("N7").Value = -6883
$ws.Range("H32").Value = 7805.2
$ws.Range("I32").Value = 2306.5
$ws.Range("K32").Value = 2306.5
$ws.Range("M32").Value = -1989.5
$ws.Range("H40").Value = 73557.14
$ws.Range("I40").Value = 201560
$ws.Range("J40").Value = 2444.4443
$ws.Range("K40").Value = 201560
$ws.Range("L40").Value = 2444.4443
$ws.Range("M40").Value = -201424
$ws.Range("N40").Value = -2716.4443
$ws.Range("H116").Value = 48480
$ws.Range("J116").Value = 48480
$ws.Range("L116").Value = 48480
$ws.Range("N116").Value = -57658
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H126").Value = 5261.3
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 6659
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 19977
$ws.Range("M126").Value = -3530
$ws.Range("N126").Value = -24917
$ws.Range("H132").Value = 1612.6184
$ws.Range("I132").Value = 1731.7164
$ws.Range("J132").Value = 726
$ws.Range("K132").Value = 5195.1492
$ws.Range("L132").Value = 2178
$ws.Range("M132").Value = -2665.1492
$ws.Range("N132").Value = -7238
$ws.Range("H136").Value = 1509.5454
$ws.Range("I136").Value = 1368.1852
$ws.Range("J136").Value = 2145.6667
$ws.Range("K136").Value = 4104.5556
$ws.Range("L136").Value = 6437.000100000001
$ws.Range("M136").Value = -1554.5556
$ws.Range("N136").Value = -11537.0001

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 19800
$ws.Range("J18").Value = 19800
$ws.Range("L18").Value = 19800
$ws.Range("N18").Value = -20146
$ws.Range("H107").Value = 59587
$ws.Range("I107").Value = 523.75
$ws.Range("J107").Value = 77760.30499999999
$ws.Range("K107").Value = 1571.25
$ws.Range("L107").Value = 233280.915
$ws.Range("M107").Value = 348.75
$ws.Range("N107").Value = -237120.915
$ws.Range("H136").Value = 941.1429000000001
$ws.Range("I136").Value = 701.2963
$ws.Range("J136").Value = 1750.625
$ws.Range("K136").Value = 2103.8889
$ws.Range("L136").Value = 5251.875
$ws.Range("M136").Value = 446.1111000000001
$ws.Range("N136").Value = -10351.875
